# cryptos.xlsx refresh — "Updated symbol list on Mon Dec 12 09:02:37 UTC 2022
# with GitHub Actions"
#
# All data cells on Sheet1 (B:G, rows 2-51) are stored as text (inlineStr),
# including numeric-looking values such as Price (column D) and Hora
# (column G). A plain `.Value = "20.80"` assignment would let Excel
# "smart type" the numeric-looking string into a real number, silently
# dropping significant trailing zeros (e.g. "20.80" -> 20.8, "0.0001500"
# -> 0.00015) and would not match the source data. To avoid that we set
# NumberFormat to Text ("@") before writing those numeric-looking values,
# then restore the cell style to "Normal" afterwards so no stray number
# format / style index is left behind on cells that otherwise carry the
# default (unstyled) format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column G ("Hora"): every data row (2-51) flips from "8" to "9". ---
$ws.Range("G2:G51").NumberFormat = "@"
$ws.Range("G2:G51").Value = "9"
$ws.Range("G2:G51").Style = "Normal"

# --- Column D ("Price"): refreshed quote per coin (only rows that moved). ---
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "281.55"
$ws.Range("D3").Value = "20.80"
$ws.Range("D4").Value = "6.230"
$ws.Range("D5").Value = "0.06141"
$ws.Range("D7").Value = "6.559"
$ws.Range("D8").Value = "1.466"
$ws.Range("D9").Value = "0.8170"
$ws.Range("D11").Value = "0.1624"
$ws.Range("D12").Value = "0.08272"
$ws.Range("D14").Value = "0.03205"
$ws.Range("D15").Value = "0.09140"
$ws.Range("D16").Value = "3.730"
$ws.Range("D17").Value = "0.001638"
$ws.Range("D18").Value = "0.04647"
$ws.Range("D19").Value = "0.006450"
$ws.Range("D20").Value = "0.006166"
$ws.Range("D21").Value = "0.001066"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("D23").Value = "3.806"
$ws.Range("D24").Value = "2.336"
$ws.Range("D25").Value = "0.3375"
$ws.Range("D26").Value = "0.1249"
$ws.Range("D40").Value = "0.04664"
$ws.Range("D41").Value = "0.007162"
$ws.Range("D42").Value = "0.1099"
$ws.Range("D43").Value = "0.003511"
$ws.Range("D44").Value = "0.01135"
$ws.Range("D45").Value = "0.00006369"
$ws.Range("D47").Value = "1.000"
$ws.Range("D49").Value = "0.00001900"

$ws.Range("D2:D51").Style = "Normal"

# --- Rows 41-43: the coin ranked list rotated by one position. ---
# CEJI (was row 41) -> KickToken (was row 42) -> BKEXToken (was row 43) -> CEJI.
# Columns B (Coin), C (Link) and E (rank+name+symbol key) are plain text,
# so a direct .Value assignment is safe — none of these strings look numeric.

# Row 41: CEJI -> KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42: KickToken -> BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43: BKEXToken -> CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"
